# Updates cryptos list values (price + volume/1h) per the commit diff.
# Rows 7 and 8 also have Coin name / Link swapped (XRP <-> LidoStakedEther).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric to Excel need to be forced to text
# first (NumberFormat "@"), matching the source file where these are stored
# as plain text (inlineStr), then the number format is reset afterwards so
# no stray style/format change is left behind on the cell.
$forceTextCells = @("D5", "D6", "D7", "D12", "D13", "D16", "D19", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D31", "D34", "D35", "D36", "D40", "D42", "D45", "D46", "D47", "D49", "D50")
foreach ($cell in $forceTextCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "70.393.63"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "3.607.92"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("D5").Value = "584.24"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("D6").Value = "191.01"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "0.632"
$ws.Range("E7").Value = "  -1.92%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.599.81"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  +2.91%  "
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "56.20"
$ws.Range("E12").Value = "  -4.17%  "
$ws.Range("D13").Value = "0.0000313"
$ws.Range("E13").Value = "  +8.27%  "
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").Value = "4.187.03"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").Value = "20.01"
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("D17").Value = "3.608.13"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "70.361.26"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").Value = "12.74"
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("D22").Value = "491.60"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "20.12"
$ws.Range("E23").Value = "  +5.74%  "
$ws.Range("D24").Value = "4.95"
$ws.Range("E24").Value = "  -7.66%  "
$ws.Range("D25").Value = "97.75"
$ws.Range("E25").Value = "  +7.50%  "
$ws.Range("D26").Value = "4.39"
$ws.Range("E26").Value = "  -1.73%  "
$ws.Range("D27").Value = "2.99"
$ws.Range("E27").Value = "  -4.14%  "
$ws.Range("D28").Value = "11.12"
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("D29").Value = "9.52"
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("D31").Value = "7.64"
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("D34").Value = "66.40"
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("D35").Value = "579.65"
$ws.Range("E35").Value = "  -9.12%  "
$ws.Range("D36").Value = "39.07"
$ws.Range("E36").Value = "  +0.98%  "
$ws.Range("D37").Value = "0.0₃0820"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("E39").Value = "  -1.35%  "
$ws.Range("D40").Value = "3.29"
$ws.Range("E40").Value = "  +20.87%  "
$ws.Range("E41").Value = "  +6.08%  "
$ws.Range("D42").Value = "3.47"
$ws.Range("E42").Value = "  -2.41%  "
$ws.Range("E43").Value = "  -6.47%  "
$ws.Range("D44").Value = "3.226.26"
$ws.Range("E44").Value = "  -2.28%  "
$ws.Range("D45").Value = "3.07"
$ws.Range("E45").Value = "  -2.27%  "
$ws.Range("D46").Value = "0.0447"
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("D47").Value = "9.65"
$ws.Range("E47").Value = "  +6.28%  "
$ws.Range("E48").Value = "  +3.12%  "
$ws.Range("D49").Value = "0.138"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("E51").Value = "  -2.68%  "

# Reset number format back to the default "Normal" style for the forced-text
# cells so only the value (not formatting) differs from the original file.
foreach ($cell in $forceTextCells) {
    $ws.Range($cell).Style = "Normal"
}

